$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all data rows (2-250)
# from serial date 45178 (2023-09-09) to 45179 (2023-09-10).
$rng = $ws.Range("C2:C250")
$rng.Value2 = 45179
